$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Kalibrasyonlar.",
    "VR",
    "VL",
    "VB",
    "IR",
    "IB",
    "VR",
    "VS",
    "VT",
    "frq"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Columns.Item(1).ColumnWidth = 18.55

$ws.Range("A14").Select()
